# Apply "normalization correction over dates" update to Sheet1.
# Rows 2-9 (columns B:F) are rescaled to new normalized values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 219.1156005859375;  "C2" = 0.1644;             "D2" = 0.1999000012874603; "E2" = 0.5945000052452087; "F2" = -0.1299999952316284
    "B3" = 171.1233978271484;  "C3" = 0.1369;             "D3" = 0.0968;             "E3" = 0.8700000047683716; "F3" = -0.1299999952316284
    "B4" = 57.83089828491211;  "C4" = 0.06759999999999999;"D4" = 0.0414;             "E4" = 0.7114999890327454; "F4" = -0.1299999952316284
    "B5" = -2.889300107955933; "C5" = -0.003;             "D5" = -0.0197;            "E5" = 0.448500007390976;  "F5" = -0.1299999952316284
    "B6" = -143.1174011230469; "C6" = -0.1132;            "D6" = -0.13;              "E6" = 0.3935999870300293; "F6" = -0.1299999952316284
    "B7" = -114.9227981567383; "C7" = -0.1157;            "D7" = -0.1299999952316284;"E7" = 0.1771000027656555; "F7" = -0.1299999952316284
    "B8" = 160.9073944091797;  "C8" = 0.1442;             "D8" = 0.1326;             "E8" = 0.4214999973773956; "F8" = -0.1299999952316284
    "B9" = 348.0477905273438;  "C9" = 0.0454;                                        "E9" = 0.8700000047683716; "F9" = -0.1299999952316284
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
